$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7
$ws.Range("B7").Value = 0.5333333333333333
$ws.Range("C7").Value = 0.8
$ws.Range("D7").Value = 0.64

# Row 8
$ws.Range("B8").Value = 0.7777777777777778
$ws.Range("C8").Value = 0.5
$ws.Range("D8").Value = 0.6086956521739131

# Row 9
$ws.Range("B9").Value = 0.625
$ws.Range("C9").Value = 0.625
$ws.Range("D9").Value = 0.625
$ws.Range("E9").Value = 0.625

# Row 10
$ws.Range("B10").Value = 0.6555555555555556
$ws.Range("C10").Value = 0.65
$ws.Range("D10").Value = 0.6243478260869566

# Row 11
$ws.Range("B11").Value = 0.6759259259259259
$ws.Range("C11").Value = 0.625
$ws.Range("D11").Value = 0.6217391304347827

# Row 12
$ws.Range("B12").Value = 0.7
$ws.Range("C12").Value = 0.7
$ws.Range("D12").Value = 0.7

# Row 13
$ws.Range("B13").Value = 0.7857142857142857
$ws.Range("C13").Value = 0.7857142857142857
$ws.Range("D13").Value = 0.7857142857142857

# Row 14
$ws.Range("B14").Value = 0.75
$ws.Range("C14").Value = 0.75
$ws.Range("D14").Value = 0.75
$ws.Range("E14").Value = 0.75

# Row 15
$ws.Range("B15").Value = 0.7428571428571429
$ws.Range("C15").Value = 0.7428571428571429
$ws.Range("D15").Value = 0.7428571428571429

# Row 16
$ws.Range("B16").Value = 0.75
$ws.Range("C16").Value = 0.75
$ws.Range("D16").Value = 0.75

# Row 22
$ws.Range("B22").Value = 0.8333333333333334
$ws.Range("C22").Value = 0.5
$ws.Range("D22").Value = 0.625

# Row 23
$ws.Range("B23").Value = 0.7222222222222222
$ws.Range("C23").Value = 0.9285714285714286
$ws.Range("D23").Value = 0.8125000000000001

# Row 24
$ws.Range("B24").Value = 0.75
$ws.Range("C24").Value = 0.75
$ws.Range("D24").Value = 0.75
$ws.Range("E24").Value = 0.75

# Row 25
$ws.Range("B25").Value = 0.7777777777777778
$ws.Range("C25").Value = 0.7142857142857143
$ws.Range("D25").Value = 0.71875

# Row 26
$ws.Range("B26").Value = 0.7685185185185185
$ws.Range("C26").Value = 0.75
$ws.Range("D26").Value = 0.734375
